$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Hiver / Eté / Année" sub-header row (row 2); this shifts
# every data row up by one and renumbers the shared strings/dimension
# automatically.
$ws.Rows.Item(2).Delete()

# Create a temporary named style whose only explicit attribute is the Font
# (size 9 Arial, same font used elsewhere in the sheet). Applying it and then
# removing the named style again leaves behind a cell format (cellXfs entry)
# that only "applies" the font - not the number format - which matches the
# formatting Excel produces for the new header cells F1:K1.
$tmpStyle = $wb.Styles.Add("TmpHeaderStyle")
$tmpStyle.Font.Size = 9
$tmpStyle.Font.Name = "Arial"

# New header row (row 1): idx / idx2 / Name / Date Start / Date End / units...
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
# E1 previously held the "(m3/s)" label (with an explicit style); the new
# label should go back to the worksheet's default (unstyled) format.
$ws.Range("E1").ClearFormats()
$ws.Range("E1").Value = "Date End"

$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

$ws.Range("F1:K1").Style = "TmpHeaderStyle"

# Drop the temporary named style; the underlying cell format stays attached
# to F1:K1 (now referencing the default "Normal" cell style), leaving behind
# exactly one new entry in cellXfs.
$tmpStyle.Delete()

# Match the selection recorded in the edited workbook.
$ws.Range("A2:K2").Select()
